# attch-email.xlsx: drop the studentId/subjectId template columns, keep
# just Name/Subject/Score (3 cols instead of 5), retarget the header/value
# rows to the remaining shared strings, fix up the merged title cell and
# the stray date-format style that was left on C3, and move the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (headers) / Row 3 (template placeholders): rename the
# surviving Name/Subject/Score cells; the old Id columns are dropped below.
$ws.Range("A2").Value = "studentName"
$ws.Range("B2").Value = "subjectName"
$ws.Range("C2").Value = "score"

$ws.Range("A3").Value = "`${student.name}"
$ws.Range("B3").Value = "`${subject.name}"
$ws.Range("C3").Value = "`${score}"

# C3 used to carry a leftover "dd/mm/yyyy;@" number format (style for the
# old date-like id column); line it back up with the rest of the value
# row (A3/B3 - centered, default number format).
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the old studentId/subjectId cells (columns D & E) from the header
# and value rows entirely - they're no longer part of the template.
$ws.Range("D2:E3").Clear()

# The title row merge shrinks from A1:E1 to A1:C1 to match the narrower
# (3-column) table.
$ws.Range("A1:E1").UnMerge()
$ws.Range("A1:C1").Merge()

# Move the saved selection/active cell.
$ws.Range("D4").Select()
